$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The stimulus picture referenced in A2 is duplicated into A3 and A4
# ("added same picture 3 times for test reasons").
$ws.Range("A3").Value = "Stimuli/083.jpg"
$ws.Range("A4").Value = "Stimuli/083.jpg"

$ws.Range("A18").Select()
